# Refactor synthetic array /3 (publipostage): update the "statut" emoji
# markers and their label from noir/rouge/orange(black/red/orange)
# to bleu/rouge/orange (blue/red/orange) colour scheme.
#
#   ⬛ -> 📘   (column A, "statut")
#   🟥 -> 📕   (column A, "statut")
#   🟧 -> 📙   (column A, "statut")
#   noir -> bleu (column B, "statut_label")

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$firstRow = $used.Row
$lastRow  = $used.Row + $used.Rows.Count - 1
$firstCol = $used.Column
$lastCol  = $used.Column + $used.Columns.Count - 1

# Locate the "statut" and "statut_label" columns from the header row so the
# edit is resilient to the exact column layout.
$statutCol = 0
$labelCol  = 0
for ($c = $firstCol; $c -le $lastCol; $c++) {
    $header = $ws.Cells.Item($firstRow, $c).Text
    if ($header -eq "statut") { $statutCol = $c }
    if ($header -eq "statut_label") { $labelCol = $c }
}
if ($statutCol -eq 0) { $statutCol = 1 }
if ($labelCol -eq 0) { $labelCol = 2 }

for ($r = $firstRow + 1; $r -le $lastRow; $r++) {
    $statutCell = $ws.Cells.Item($r, $statutCol)   # "statut" column
    $labelCell  = $ws.Cells.Item($r, $labelCol)    # "statut_label" column

    $statutText = $statutCell.Text
    if ($statutText -eq "⬛") {
        $statutCell.Value = "📘"
    } elseif ($statutText -eq "🟥") {
        $statutCell.Value = "📕"
    } elseif ($statutText -eq "🟧") {
        $statutCell.Value = "📙"
    }

    if ($labelCell.Text -eq "noir") {
        $labelCell.Value = "bleu"
    }
}
